$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.337.77"
$ws.Range("E2").Value = "  -3.28%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.498.51"
$ws.Range("E3").Value = "  -4.77%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "607.34"
$ws.Range("E5").Value = "  -2.43%  "

# Row 6 - Solana
$ws.Range("D6").Value = "149.25"
$ws.Range("E6").Value = "  -6.38%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.497.08"
$ws.Range("E7").Value = "  -4.72%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -3.08%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.81%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  -2.74%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -4.54%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.088.13"
$ws.Range("E14").Value = "  -4.76%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "31.48"
$ws.Range("E15").Value = "  -2.70%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.492.34"
$ws.Range("E16").Value = "  -4.68%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.254.97"
$ws.Range("E17").Value = "  -3.43%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.96%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "15.02"
$ws.Range("E20").Value = "  -5.38%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "446.16"
$ws.Range("E21").Value = "  -5.17%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -12.65%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -4.30%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "77.26"
$ws.Range("E24").Value = "  -3.03%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +4.56%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.12%  "

# Row 27 - WrappedeETH
$ws.Range("D27").Value = "3.637.04"
$ws.Range("E27").Value = "  -4.76%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "10.16"
$ws.Range("E28").Value = "  -8.07%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -4.75%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -6.26%  "

# Row 33 - Kaspa
$ws.Range("D33").Value = "0.163"
$ws.Range("E33").Value = "  -0.26%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "25.69"
$ws.Range("E34").Value = "  -3.35%  "

# Row 35 - NEARProtocol
$ws.Range("D35").Value = "6.14"
$ws.Range("E35").Value = "  -4.06%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "1.86"
$ws.Range("E36").Value = "  -6.36%  "

# Row 37 - RenzoRestakedETH
$ws.Range("D37").Value = "3.485.95"
$ws.Range("E37").Value = "  -5.11%  "

# Row 38 - Aptos
$ws.Range("D38").Value = "8.00"
$ws.Range("E38").Value = "  -3.36%  "

# Row 39 - USDe
$ws.Range("E39").Value = "  +0.11%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  +0.01%  "

# Rows 41-42: Monero/Stacks swap position with updated values
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "2.19"
$ws.Range("E41").Value = "  -1.42%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "173.27"
$ws.Range("E42").Value = "  -2.67%  "

# Row 43 - Hedera
$ws.Range("D43").Value = "0.0876"
$ws.Range("E43").Value = "  -1.70%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "5.43"
$ws.Range("E44").Value = "  -6.05%  "

# Row 45 - Mantle
$ws.Range("E45").Value = "  -4.68%  "

# Row 46 - OKB
$ws.Range("D46").Value = "45.45"
$ws.Range("E46").Value = "  -2.67%  "

# Row 47 - ONDO
$ws.Range("E47").Value = "  +5.26%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "26.95"
$ws.Range("E48").Value = "  -6.32%  "

# Row 49 - dogwifhat
$ws.Range("E49").Value = "  -5.01%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  -4.05%  "

# Row 51 - SuiNetwork
$ws.Range("E51").Value = "  -3.19%  "
